$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Finished presentation example: swap the old "Incorta" google-search
# example for the new "Selenium" one, and drop the now-unused "Menu" demo
# columns (D/E) that belonged to an earlier, unfinished example ---

# Row 1 headers
$ws.Range("A1").Value = "Browser"
$ws.Range("B1").Value = "Search Data"
$ws.Range("C1").Value = "Texts"

# Row 2 - the actual search result data
$ws.Range("A2").Value = "Chrome"
$ws.Range("B2").Value = "Selenium"
$ws.Range("C2").Value = "What is Selenium? Introduction to Selenium Automation Testing"

# Row 3
$ws.Range("A3").Value = "Firefox"

# The old "Menu" / "TEAM" / "ABOUT" / "Menu Item" example columns are no
# longer needed now that the presentation examples are finished
$null = $ws.Range("D1:E2").ClearContents()

# Widen column C to fit the new, longer text and leave the selection on
# the cell that was last edited
$ws.Columns.Item(3).ColumnWidth = 58
$null = $ws.Range("C2").Select()

Write-Output "edit complete"
